$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply new custom date format (dd-mmm-yyyy) to the Invoice Date column rows 2-6
$ws.Range("D2:D6").NumberFormat = "dd\-mmm\-yyyy"

# Row 4 data correction: Invoice No and Order Tracking No now match row 2's values
$ws.Range("C4").Value = 123456452
$ws.Range("E4").Value = "qw1239967"

# Update the active selection shown when the sheet was saved
$ws.Range("D9").Select()
